$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new variable entry to the list, below the existing "for D: ..." row.
# Copy B38's formatting (fill/style) down to B39 first, then overwrite its
# value/text with the new variable description.
$ws.Range("B38").Copy()
$ws.Range("B39").PasteSpecial()
$ws.Range("B39").Value = "historical salary of each player"

# Clear clipboard / marching-ants state left over from the copy above.
$excel.CutCopyMode = $false

# Match the author's resulting selection (one row below the new entry).
$ws.Range("B40").Select()
